$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Degree")

$ws.Range("F2").Formula = "=D2&`$A`$1&""=""&A2&"",""&`$B`$1&""='""&B2&""','""&`$C`$1&""='""&C2&""'),"""
$ws.Range("F3").Formula = "=D3&`$A`$1&""=""&A3&"",""&`$B`$1&""='""&B3&""','""&`$C`$1&""='""&C3&""'),"""
$ws.Range("F3").Copy()
$ws.Range("F4:F5").PasteSpecial()

Write-Host "done"
